$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'329.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.71%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'40.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'10.97%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.974"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'17.33%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.08136"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.64%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'4.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.02%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'8.764"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.52%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'1.969"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7.19%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'-0.21%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.9488"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.00%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.1318"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'16.13%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.1993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.52%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09374"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.53%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03479"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.35%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09609"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001311"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.79%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006540"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'11.47%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.357"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.03%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3542"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.25%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'9.557"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'52.15%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'9.60%"
$ws.Range("E21").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04442"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.59%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.001261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'4.90%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.004441"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.85%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.0001095"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-18.17%"
$ws.Range("E26").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.02479"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'17.12%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.05294"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.24%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007524"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.20%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1437"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.49%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.009047"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.43%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.002059"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.56%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.01032"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'33.07%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00006841"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'7.45%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.003509"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'21.76%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'6.73%"
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'-0.22%"
$ws.Range("E51").Style = "Normal"
